# "Modified DSL for EB"
#
# The EB_Notification_JS test-DSL spreadsheet referenced the old Rhomobile
# Compliance app id / notification target ("Compliance") in a handful of
# test-step cells. This updates those DSL strings to the new Enterprise
# Browser app id / notification target, and clears the stray "Pass"/"Fail"
# values that were left over in column J (Results).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Update-CellText {
    param($Address, $Find, $Replace)
    $current = $ws.Range($Address).Text
    $ws.Range($Address).Value = $current.Replace($Find, $Replace)
}

# VT200_0683 (row 5) and VT200_0687 (row 9) launch the app via its package
# name - point them at the new Enterprise Browser package.
Update-CellText "G5" "com.rhomobile.compliancetest_js" "com.symbol.enterprisebrowser"
Update-CellText "G9" "com.rhomobile.compliancetest_js" "com.symbol.enterprisebrowser"

# VT200_0685 (row 7) and VT200_0686 (row 8) swipe-and-tap a notification by
# its title - the notification is now titled "EnterpriseBrowser" instead of
# "Compliance".
Update-CellText "G7" "SwipeNotificationAndTap(Compliance)" "SwipeNotificationAndTap(EnterpriseBrowser)"
Update-CellText "G8" "SwipeNotificationAndTap(Compliance)" "SwipeNotificationAndTap(EnterpriseBrowser)"

# Column J ("Results") held stray Pass/Fail markers that are no longer
# used - clear them out.
$ws.Range("J2").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("J6").Value = ""
$ws.Range("J7").Value = ""
$ws.Range("J8").Value = ""
$ws.Range("J9").Value = ""
$ws.Range("J11").Value = ""

# Move the saved selection back to the top of the sheet.
$ws.Range("G1").Select()
